# Update "想去人数" (F column) counts across all four sheets to match
# the regenerated gh-pages output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 36
$ws.Range("F5").Value = 1936
$ws.Range("F6").Value = 1298
$ws.Range("F7").Value = 1298
$ws.Range("F13").Value = 1706
$ws.Range("F15").Value = 1856
$ws.Range("F17").Value = 1028
$ws.Range("F18").Value = 48
$ws.Range("F19").Value = 514
$ws.Range("F20").Value = 1598
$ws.Range("F22").Value = 22
$ws.Range("F25").Value = 2383
$ws.Range("F26").Value = 434
$ws.Range("F28").Value = 1018
$ws.Range("F29").Value = 4542
$ws.Range("F30").Value = 104
$ws.Range("F35").Value = 1241

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F20").Value = 185
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 13

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F9").Value = 3090
$ws.Range("F10").Value = 607
$ws.Range("F14").Value = 43

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 36
$ws.Range("F8").Value = 3090
$ws.Range("F9").Value = 607
$ws.Range("F11").Value = 1936
$ws.Range("F14").Value = 43
$ws.Range("F15").Value = 1298
$ws.Range("F21").Value = 1706
$ws.Range("F23").Value = 1856
$ws.Range("F24").Value = 1028
$ws.Range("F25").Value = 48
$ws.Range("F26").Value = 514
$ws.Range("F28").Value = 1598
$ws.Range("F30").Value = 185
$ws.Range("F31").Value = 13
$ws.Range("F32").Value = 22
$ws.Range("F36").Value = 2383
$ws.Range("F37").Value = 434
$ws.Range("F41").Value = 4542
$ws.Range("F51").Value = 1241
